# Clear out the (previously entered) production-defect rows 286:307.
# Values/formulas are removed in columns B-P, but the existing cell
# styles (s="...") are left intact since ClearContents() only wipes
# content, not formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B286:P307").ClearContents()

# Move/restore the sheet's active selection to D297 (was F306).
[void]$ws.Range("D297").Select()
